# Rerun snapshot with 2025 commercial indicators.
# Updates the longfin squid indicator table: bumps the "long term" baseline
# years from 1996-2024 to 1996-2025, refreshes the regenerated chart image
# filenames, and rewrites the commercial-vessel and commercial-revenue
# narrative text/labels for the new 2025 data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Commercial landings
$ws.Range("B2").Value = "Near long term (1996-2025) average"
$ws.Range("D2").Value = "Commercial_LONGFINSQUID_Landings_LBS_2026-02-25.png"

# Row 3 - Number of commercial vessels
$ws.Range("B3").Value = "Below long term (1996-2025) average"
$ws.Range("C3").Value = "Number of commercial vessels has been steadily`ndecreasing since around 2000 (although has slightly increased in 2025) consistent with`ndecreasing fleet diversity and continued risk to`nfishery resilience [7]. Permit requalification in 2019 and a decrease in the post-closure trip  limit for trimester 2 may cap participation, although these actions were designed to accommodate recent fishing trends and activity. "
$ws.Range("D3").Value = "N_Commercial_Vessels_Landing_LONGFINSQUID_2026-02-25.png"

# Row 4 - Commercial revenue
$ws.Range("A4").Value = "Commercial revenue (millions, inflation adjusted to 2025 USD)"
$ws.Range("B4").Value = "Below long term (1996-2025) average"
$ws.Range("C4").Value = "Average longfin ex-vessel prices in 2024 increased slightly from 2023 (+10%). Commercial revenue has decreased since 2022, driven by the overall`ndecrease in landings by 23% [7]."
$ws.Range("D4").Value = "TOTALANNUALREV_LONGFINSQUID_2025Dols_2026-02-25.png"

# Row heights adjusted to fit the rewrapped text in rows 3 and 4
$ws.Rows.Item(3).RowHeight = 185.4
$ws.Rows.Item(4).RowHeight = 79.8

# Selection moved to D4
$ws.Range("D4").Select()
